$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: "Group and Self Assessment"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Group and Self Assessment")

# Grade corrections (4 -> 5) on the self-assessment diagonal
$ws1.Range("D10").Value = 5
$ws1.Range("E11").Value = 5
$ws1.Range("F12").Value = 5
$ws1.Range("G13").Value = 5

# View state: scroll position + selection
$ws1.Range("H13").Select()
$ws1.Application.ActiveWindow.ScrollRow = 8

# ---------------------------------------------------------------------------
# Sheet: "User Stories"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("User Stories")

$ws2.Range("B22").Value = 1201925
$ws2.Range("C22").Value = 5
$ws2.Range("B23").Value = 1222183
$ws2.Range("C23").Value = 5
$ws2.Range("B24").Value = 1230420
$ws2.Range("C24").Value = 5

$ws2.Range("C26").Value = 5
$ws2.Range("C28").Value = 5
$ws2.Range("C29").Value = 5
$ws2.Range("C30").Value = 5
$ws2.Range("C31").Value = 5
$ws2.Range("C32").Value = 5
$ws2.Range("C33").Value = 5
$ws2.Range("C34").Value = 5

# View state: scroll position + selection
$ws2.Range("B24").Select()
$ws2.Application.ActiveWindow.ScrollRow = 20

# ---------------------------------------------------------------------------
# Sheet: "Project Development"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Project Development")

$ws3.Range("D5").Value = 5
$ws3.Range("C6").Value = 5
$ws3.Range("D6").Value = 5
$ws3.Range("E6").Value = 5
$ws3.Range("F6").Value = 5

# View state: selection (topLeftCell resets to default/top)
$ws3.Range("H8").Select()

# ---------------------------------------------------------------------------
# Activate "User Stories" tab last so it becomes the workbook's active sheet
# ---------------------------------------------------------------------------
$ws2.Activate()
